$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H108").Value = 31197.334
$ws.Range("J108").Value = 31197.334
$ws.Range("L108").Value = 31197.334
$ws.Range("N108").Value = -38877.334
$ws.Range("H124").Value = 53315.332
$ws.Range("J124").Value = 53315.332
$ws.Range("L124").Value = 53315.332
$ws.Range("N124").Value = -63135.332
$ws.Range("H126").Value = 47772
$ws.Range("J126").Value = 47772
$ws.Range("L126").Value = 47772
$ws.Range("N126").Value = -57652
$ws.Range("H128").Value = 57484
$ws.Range("J128").Value = 57484
$ws.Range("L128").Value = 57484
$ws.Range("N128").Value = -67444
$ws.Range("H130").Value = 54992
$ws.Range("J130").Value = 54992
$ws.Range("L130").Value = 54992
$ws.Range("N130").Value = -65032

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H111").Value = 47888
$ws.Range("J111").Value = 47888
$ws.Range("L111").Value = 47888
$ws.Range("N111").Value = -56068
$ws.Range("H117").Value = 46998
$ws.Range("J117").Value = 46998
$ws.Range("L117").Value = 46998
$ws.Range("N117").Value = -56176

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H126").Value = 50776
$ws.Range("J126").Value = 50776
$ws.Range("L126").Value = 50776
$ws.Range("N126").Value = -60656

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 49776
$ws.Range("J20").Value = 49776
$ws.Range("L20").Value = 49776
$ws.Range("N20").Value = -50248
$ws.Range("H30").Value = 49776
$ws.Range("J30").Value = 49776
$ws.Range("L30").Value = 49776
$ws.Range("N30").Value = -49958
$ws.Range("H75").Value = 19999.889
$ws.Range("J75").Value = 19999.889
$ws.Range("L75").Value = 19999.889
$ws.Range("N75").Value = -21995.889
$ws.Range("H78").Value = 19999.889
$ws.Range("J78").Value = 19999.889
$ws.Range("L78").Value = 59999.667
$ws.Range("N78").Value = -69983.667
$ws.Range("H116").Value = 44974
$ws.Range("J116").Value = 44974
$ws.Range("L116").Value = 44974
$ws.Range("N116").Value = -54152
$ws.Range("H128").Value = 49776
$ws.Range("J128").Value = 49776
$ws.Range("L128").Value = 49776
$ws.Range("N128").Value = -59736
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H133").Value = 29630.666
$ws.Range("J133").Value = 29630.666
$ws.Range("L133").Value = 29630.666
$ws.Range("N133").Value = -34690.666
$ws.Range("H137").Value = 61923.8
$ws.Range("J137").Value = 61923.8
$ws.Range("L137").Value = 61923.8
$ws.Range("N137").Value = -72123.8
$ws.Range("H138").Value = 45657.8
$ws.Range("J138").Value = 45657.8
$ws.Range("L138").Value = 45657.8
$ws.Range("N138").Value = -55937.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5177.592
$ws.Range("I70").Value = 5264.7646
$ws.Range("J70").Value = 4980
$ws.Range("K70").Value = 5264.7646
$ws.Range("L70").Value = 4980
$ws.Range("M70").Value = -4994.7646
$ws.Range("N70").Value = -5520
$ws.Range("H73").Value = 5177.592
$ws.Range("I73").Value = 5264.7646
$ws.Range("J73").Value = 4980
$ws.Range("K73").Value = 5264.7646
$ws.Range("L73").Value = 4980
$ws.Range("M73").Value = -4328.7646
$ws.Range("N73").Value = -6852
$ws.Range("H110").Value = 48702
$ws.Range("J110").Value = 48702
$ws.Range("L110").Value = 48702
$ws.Range("N110").Value = -56882
$ws.Range("H130").Value = 47692
$ws.Range("J130").Value = 47692
$ws.Range("L130").Value = 47692
$ws.Range("N130").Value = -57732
$ws.Range("H137").Value = 59799.5
$ws.Range("J137").Value = 59799.5
$ws.Range("L137").Value = 59799.5
$ws.Range("N137").Value = -69999.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2471.7917
$ws.Range("I16").Value = 1714.0435
$ws.Range("J16").Value = 19900
$ws.Range("K16").Value = 1714.0435
$ws.Range("L16").Value = 19900
$ws.Range("M16").Value = -1544.0435
$ws.Range("N16").Value = -20240
$ws.Range("H22").Value = 416.86667
$ws.Range("I22").Value = 406.66666
$ws.Range("J22").Value = 457.66666
$ws.Range("K22").Value = 406.66666
$ws.Range("L22").Value = 457.66666
$ws.Range("M22").Value = -111.66666
$ws.Range("N22").Value = -1047.66666
$ws.Range("H27").Value = 416.86667
$ws.Range("I27").Value = 406.66666
$ws.Range("J27").Value = 457.66666
$ws.Range("K27").Value = 406.66666
$ws.Range("L27").Value = 457.66666
$ws.Range("M27").Value = -299.66666
$ws.Range("N27").Value = -671.66666
$ws.Range("H111").Value = 39544.75
$ws.Range("J111").Value = 39544.75
$ws.Range("L111").Value = 39544.75
$ws.Range("N111").Value = -47724.75
$ws.Range("H133").Value = 35141.715
$ws.Range("J133").Value = 35141.715
$ws.Range("L133").Value = 35141.715
$ws.Range("N133").Value = -40201.715
$ws.Range("H137").Value = 48745.75
$ws.Range("J137").Value = 48745.75
$ws.Range("L137").Value = 48745.75
$ws.Range("N137").Value = -58945.75
$ws.Range("H138").Value = 55342.855
$ws.Range("J138").Value = 55342.855
$ws.Range("L138").Value = 55342.855
$ws.Range("N138").Value = -65622.85500000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 29324
$ws.Range("J86").Value = 29324
$ws.Range("L86").Value = 29324
$ws.Range("N86").Value = -31570
$ws.Range("H89").Value = 29324
$ws.Range("J89").Value = 29324
$ws.Range("L89").Value = 146620
$ws.Range("N89").Value = -157852
$ws.Range("H108").Value = 48292.668
$ws.Range("J108").Value = 48292.668
$ws.Range("L108").Value = 48292.668
$ws.Range("N108").Value = -55972.668
$ws.Range("H110").Value = 50172.8
$ws.Range("J110").Value = 50172.8
$ws.Range("L110").Value = 50172.8
$ws.Range("N110").Value = -58352.8
$ws.Range("H116").Value = 48696
$ws.Range("J116").Value = 48696
$ws.Range("L116").Value = 48696
$ws.Range("N116").Value = -57874
$ws.Range("H119").Value = 48698
$ws.Range("J119").Value = 48698
$ws.Range("L119").Value = 48698
$ws.Range("N119").Value = -58374
